$d = $word.ActiveDocument

$d.Content.Find.Execute("59×74=", $true, $false, $false, $false, $false, $true, 1, $false, "34×67=", 2) | Out-Null
$d.Content.Find.Execute("70×43=", $true, $false, $false, $false, $false, $true, 1, $false, "56×78=", 2) | Out-Null
$d.Content.Find.Execute("85×29=", $true, $false, $false, $false, $false, $true, 1, $false, "66×52=", 2) | Out-Null
$d.Content.Find.Execute("21×37=", $true, $false, $false, $false, $false, $true, 1, $false, "65×97=", 2) | Out-Null
$d.Content.Find.Execute("43×70=", $true, $false, $false, $false, $false, $true, 1, $false, "61×93=", 2) | Out-Null
$d.Content.Find.Execute("88×53=", $true, $false, $false, $false, $false, $true, 1, $false, "99×17=", 2) | Out-Null
$d.Content.Find.Execute("37×86=", $true, $false, $false, $false, $false, $true, 1, $false, "69×38=", 2) | Out-Null
$d.Content.Find.Execute("96×78=", $true, $false, $false, $false, $false, $true, 1, $false, "51×85=", 2) | Out-Null
$d.Content.Find.Execute("52×86=", $true, $false, $false, $false, $false, $true, 1, $false, "40×83=", 2) | Out-Null
$d.Content.Find.Execute("82×43=", $true, $false, $false, $false, $false, $true, 1, $false, "62×13=", 2) | Out-Null
$d.Content.Find.Execute("99×91=", $true, $false, $false, $false, $false, $true, 1, $false, "29×18=", 2) | Out-Null
$d.Content.Find.Execute("29×45=", $true, $false, $false, $false, $false, $true, 1, $false, "24×96=", 2) | Out-Null
$d.Content.Find.Execute("60×28=", $true, $false, $false, $false, $false, $true, 1, $false, "35×57=", 2) | Out-Null
$d.Content.Find.Execute("40×39=", $true, $false, $false, $false, $false, $true, 1, $false, "85×48=", 2) | Out-Null
$d.Content.Find.Execute("63×89=", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=", 2) | Out-Null
$d.Content.Find.Execute("81×88=", $true, $false, $false, $false, $false, $true, 1, $false, "34×66=", 2) | Out-Null
$d.Content.Find.Execute("64×97=", $true, $false, $false, $false, $false, $true, 1, $false, "50×63=", 2) | Out-Null
$d.Content.Find.Execute("64×98=", $true, $false, $false, $false, $false, $true, 1, $false, "41×96=", 2) | Out-Null
$d.Content.Find.Execute("39×71=", $true, $false, $false, $false, $false, $true, 1, $false, "64×67=", 2) | Out-Null
$d.Content.Find.Execute("17×59=", $true, $false, $false, $false, $false, $true, 1, $false, "61×35=", 2) | Out-Null
$d.Content.Find.Execute("37×91=", $true, $false, $false, $false, $false, $true, 1, $false, "63×54=", 2) | Out-Null
$d.Content.Find.Execute("26×86=", $true, $false, $false, $false, $false, $true, 1, $false, "24×72=", 2) | Out-Null
$d.Content.Find.Execute("79×67=", $true, $false, $false, $false, $false, $true, 1, $false, "13×24=", 2) | Out-Null
$d.Content.Find.Execute("92×17=", $true, $false, $false, $false, $false, $true, 1, $false, "62×44=", 2) | Out-Null
$d.Content.Find.Execute("13×17=", $true, $false, $false, $false, $false, $true, 1, $false, "62×80=", 2) | Out-Null
